# Add the new MQTT-related translation rows (33-38) to the "Translation"
# sheet, mirroring the rows already present (TEXT ID / TYPOGRAPHY NAME /
# ALIGNMENT / DIRECTION / GB columns = B..F).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$newline = [char]10
$invalidCred = "Invalid credencial, " + $newline + "try again please"

$rows = @(
    @{ Row = 33; Id = "SingleUseId28"; Typo = "Typography_label"; Align = "Center"; Dir = "LTR"; Text = $invalidCred },
    @{ Row = 34; Id = "SingleUseId29"; Typo = "Typography_label"; Align = "Center"; Dir = "LTR"; Text = "OK" },
    @{ Row = 35; Id = "SingleUseId30"; Typo = "Typography_label"; Align = "Center"; Dir = "LTR"; Text = "Ok" },
    @{ Row = 36; Id = "SingleUseId31"; Typo = "Typography_label"; Align = "Center"; Dir = "LTR"; Text = "Cancel" },
    @{ Row = 37; Id = "SingleUseId32"; Typo = "Typography_label"; Align = "Left";   Dir = "LTR"; Text = "You can dispatch <value> liters" },
    @{ Row = 38; Id = "SingleUseId33"; Typo = "Typography_label"; Align = "Left";   Dir = "LTR"; Text = "0" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.Id     # column B - TEXT ID
    $ws.Cells.Item($r.Row, 3).Value = $r.Typo   # column C - TYPOGRAPHY NAME
    $ws.Cells.Item($r.Row, 4).Value = $r.Align  # column D - ALIGNMENT
    $ws.Cells.Item($r.Row, 5).Value = $r.Dir    # column E - DIRECTION

    # The "GB" text column holds display text for the UI. Some of these
    # look like plain numbers (e.g. the "0" keypad label), but they must
    # stay text (same treatment as the existing "0".."9" keypad labels
    # already in this sheet), so force text storage for those.
    $cell = $ws.Cells.Item($r.Row, 6)
    if ($r.Text -match '^-?\d+(\.\d+)?$') {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $r.Text
}
